# Generate Report for Handoff
# The localization status report moved from "In Translation" to
# "Ready for handoff" and the handoff timestamps were refreshed. Update the
# Overview rollup sheet plus the per-language (zh-cn / de-de) detail sheets,
# then widen the status columns so the longer "Ready for handoff" text fits
# (mirrors the auto-fit that produced the new column widths).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps ------------------------------------------
$wsOverview.Range("G2").Value = "2016-08-27 00:58:31"
$wsZhCn.Range("H2").Value     = "2016-08-27 00:58:27"
$wsDeDe.Range("H2").Value     = "2016-08-27 00:58:31"

# --- Widen the status columns to fit "Ready for handoff" -------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33
